{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in bold + color (#2C3E50) across the resume body. For each\n// affected paragraph, the existing single run of text is split (via\n// paragraph.search) around the metric substrings and only the metric\n// substrings get the bold + color formatting applied; the surrounding plain\n// text is left untouched other than being split into separate runs.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Map of paragraph index (within context.document.body.paragraphs) ->\n// ordered list of literal substrings to bold + color within that paragraph.\n// Order matters only in that every listed term must be unique/unambiguous\n// within its paragraph (verified against the source document).\nconst targets = [\n  // \"\u2022 Discovered systematic race coding errors ... from 23% to 64%\"\n  { index: 9, terms: [\"23%\", \"64%\"] },\n  // \"\u2022 Utilized advanced sampling methods ... from \u00b14.2% to \u00b12.1%, increasing\n  //   voter turnout prediction accuracy from 71% to 87%, ...\"\n  { index: 11, terms: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"] },\n  // \"\u2022 Trigonometric algorithm ... reduced mapping costs by 73.5%, saving\n  //   campaigns and organizations $4.7M ...\"\n  { index: 12, terms: [\"73.5%\", \"$4.7M\"] },\n  // \"\u2022 Built real-time FEC analysis systems ... valued over $2 trillion\"\n  { index: 13, terms: [\"$2\"] },\n  // \"\u2022 Modernized legacy ETL processes ... reducing processing time by 57%\"\n  { index: 18, terms: [\"57%\"] },\n  // \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts ...\"\n  { index: 49, terms: [\"12,847\"] },\n  // \"\u2022 Revenue generation: Delivered $4.9M additional revenue ...\"\n  { index: 51, terms: [\"$4.9M\"] },\n  // \"\u2022 23% conversion rate improvement\"\n  { index: 52, terms: [\"23%\"] },\n];\n\nconst BOLD_COLOR = \"#2C3E50\";\n\nfor (const { index, terms } of targets) {\n  const paragraph = paragraphs.items[index];\n  for (const term of terms) {\n    const found = paragraph.search(term, { matchCase: true, matchWholeWord: false });\n    found.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < found.items.length; i++) {\n      const range = found.items[i];\n      range.font.bold = true;\n      range.font.color = BOLD_COLOR;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in bold + color (#2C3E50) across the resume body. For each\n# affected paragraph, we Find() the literal metric substring within that\n# paragraph's Range and apply Bold + Font.Color directly to the found\n# (now-collapsed-to-match) Range. Word automatically splits the paragraph's\n# run(s) so only the matched substring carries the new run-level formatting,\n# which is exactly the run-splitting behavior captured in the diff.\n\nfunction RGBColor($r, $g, $b) {\n    # Word's Font.Color (and the underlying w:color value) is documented in\n    # RGB hex (e.g. 2C3E50) but the OLE/COM long value is packed as BGR.\n    return $r + ($g * 256) + ($b * 65536)\n}\n\n$d = $word.ActiveDocument\n$highlightColor = RGBColor 0x2C 0x3E 0x50\n\n# Paragraphs collection is 1-indexed. For each affected paragraph, list the\n# literal metric substrings (in left-to-right order) that must become bold +\n# colored. Every term occurs exactly once in its paragraph, so a fresh\n# Find.Execute per term (scoped to that paragraph's Range) unambiguously\n# locates it regardless of edits already made earlier in the same paragraph.\n$targets = @(\n    # \"Discovered systematic race coding errors ... from 23% to 64%\"\n    @{ Para = 10; Terms = @(\"23%\", \"64%\") },\n    # \"Utilized advanced sampling methods ... from \u00b14.2% to \u00b12.1%, increasing\n    #  voter turnout prediction accuracy from 71% to 87%, ...\"\n    @{ Para = 12; Terms = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\") },\n    # \"Trigonometric algorithm ... reduced mapping costs by 73.5%, saving\n    #  campaigns and organizations $4.7M ...\"\n    @{ Para = 13; Terms = @(\"73.5%\", \"$4.7M\") },\n    # \"Built real-time FEC analysis systems ... valued over $2 trillion\"\n    @{ Para = 14; Terms = @(\"$2\") },\n    # \"Modernized legacy ETL processes ... reducing processing time by 57%\"\n    @{ Para = 19; Terms = @(\"57%\") },\n    # \"Platform impact: Built redistricting system serving 12,847 analysts ...\"\n    @{ Para = 50; Terms = @(\"12,847\") },\n    # \"Revenue generation: Delivered $4.9M additional revenue ...\"\n    @{ Para = 52; Terms = @(\"$4.9M\") },\n    # \"23% conversion rate improvement\"\n    @{ Para = 53; Terms = @(\"23%\") }\n)\n\nforeach ($t in $targets) {\n    $paraIndex = $t.Para\n    foreach ($term in $t.Terms) {\n        $rng = $d.Paragraphs($paraIndex).Range\n        $found = $rng.Find.Execute($term)\n        if ($found) {\n            $rng.Bold = $true\n            $rng.Font.Color = $highlightColor\n        }\n    }\n}\n"}
